# CORE_holdings.xlsx — refresh the model-holdings snapshot:
#  - bump the "as of" date in the confidential disclaimer from 2021-05-17 to 2021-05-18
#  - update the Weight (D) and Percent Change (E) figures for each row

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet ships protected; temporarily unprotect so the cells can be written.
$ws.Unprotect()

# --- Disclaimer date -------------------------------------------------------
# Surgically replace just the date substring inside the multi-line disclaimer
# cell rather than rewriting the whole string.
$disclaimer = $ws.Range("A11").Value2
$oldDate = "2021-05-17"
$newDate = "2021-05-18"
$idx = $disclaimer.IndexOf($oldDate)
if ($idx -ge 0) {
    $ws.Range("A11").Characters($idx + 1, $oldDate.Length).Text = $newDate
}

# --- Weight / Percent Change updates ---------------------------------------
$ws.Range("D2").Value = 0.5019744104819298
$ws.Range("E2").Value = -0.008549292011755361

$ws.Range("D3").Value = 0.2418844766512078
$ws.Range("E3").Value = -0.007995262066923359

$ws.Range("D4").Value = 0.09489457396681343
$ws.Range("E4").Value = -0.007967623624636322

$ws.Range("D5").Value = 0.1041639877967354
$ws.Range("E5").Value = -0.01316752011704458

$ws.Range("D6").Value = 0.03029194681693103
$ws.Range("E6").Value = -0.01273464767474786

$ws.Range("D7").Value = 0.02679060428638238
$ws.Range("E7").Value = -0.00789988267500974

$ws.Range("D8").Value = 0.9999999999999998
$ws.Range("E8").Value = -0.008950521153309521

# Restore (best-effort) sheet protection so the workbook isn't left open.
$ws.Protect()
